{"js": "// Update the date line and the division problems in the table, in place,\n// replacing each text run's content without altering the document\n// structure (same number of paragraphs/rows/cells, same run formatting).\n\nconst body = context.document.body;\n\n// --- 1. Update the date paragraph -----------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\nif (dateParagraph.text.trim() === \"2023-09-13 Wednesday\") {\n  dateParagraph.insertText(\"2023-09-14 Thursday\", \"Replace\");\n}\n\n// --- 2. Update the division-problem table ----------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row-major grid of the new cell values. Only rows 0, 4, 8, 12, 16 hold\n// text (the other rows are blank spacer rows) and every row has 5 cells.\nconst newValues = {\n  0: [\"59\u00f74=\", \"24\u00f73=\", \"10\u00f74=\", \"42\u00f77=\", \"56\u00f78=\"],\n  4: [\"71\u00f76=\", \"98\u00f72=\", \"91\u00f79=\", \"15\u00f72=\", \"75\u00f79=\"],\n  8: [\"30\u00f73=\", \"11\u00f79=\", \"13\u00f77=\", \"47\u00f73=\", \"96\u00f75=\"],\n  12: [\"84\u00f73=\", \"43\u00f75=\", \"69\u00f74=\", \"67\u00f75=\", \"65\u00f78=\"],\n  16: [\"49\u00f75=\", \"16\u00f72=\", \"92\u00f78=\", \"40\u00f72=\", \"20\u00f76=\"],\n};\n\nfor (const rowIndex of Object.keys(newValues)) {\n  const r = Number(rowIndex);\n  const rowValues = newValues[rowIndex];\n  for (let c = 0; c < rowValues.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = rowValues[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the division problems in the table, in place,\n# keeping the same document/table structure and run formatting.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the date paragraph ------------------------------------\n$dateRange = $d.Paragraphs.Item(1).Range\n$dateRange.MoveEnd(1, -1) | Out-Null\n$dateRange.Text = \"2023-09-14 Thursday\"\n\n# --- 2. Update the division-problem table -----------------------------\n$table = $d.Tables.Item(1)\n\n# Row-major grid of the new cell values (1-based row/col indices).\n# Only rows 1, 5, 9, 13, 17 hold text (the rest are blank spacer rows),\n# and each of those rows has 5 cells.\n$newValues = @{\n    1  = @(\"59\u00f74=\", \"24\u00f73=\", \"10\u00f74=\", \"42\u00f77=\", \"56\u00f78=\")\n    5  = @(\"71\u00f76=\", \"98\u00f72=\", \"91\u00f79=\", \"15\u00f72=\", \"75\u00f79=\")\n    9  = @(\"30\u00f73=\", \"11\u00f79=\", \"13\u00f77=\", \"47\u00f73=\", \"96\u00f75=\")\n    13 = @(\"84\u00f73=\", \"43\u00f75=\", \"69\u00f74=\", \"67\u00f75=\", \"65\u00f78=\")\n    17 = @(\"49\u00f75=\", \"16\u00f72=\", \"92\u00f78=\", \"40\u00f72=\", \"20\u00f76=\")\n}\n\nforeach ($rowIndex in $newValues.Keys) {\n    $rowValues = $newValues[$rowIndex]\n    for ($c = 0; $c -lt $rowValues.Length; $c++) {\n        $table.Cell($rowIndex, $c + 1).Range.Text = $rowValues[$c]\n    }\n}\n"}
